$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete worksheet row 768 (the "二回失敗しているペンギン" post entry),
# shifting all subsequent rows up by one.
$ws.Rows.Item(768).Delete()
